# "equine cv changed 3.0"
# Applies the CV edits to slide 1 (Varvara Lazarenko PhD Aneuploidy CV):
#  - reflow / resize the education text box and the two adjacent date labels
#  - expand the "Skills" list in the education text box
#  - tweak the wording of the personal statement

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Shape "Google Shape;99;p1": big education/skills text box ---
$eduShape = $s.Shapes.Item("Google Shape;99;p1")
$eduShape.Left = 113.0728
$eduShape.Width = 479.2904
$eduShape.Height = 301.7091

$eduText = $eduShape.TextFrame.TextRange
$oldSkills = ": PCR, gel electrophoresis, wire myography, western blotting, ELISA, immunohistochemistry, intracellular recording (microelectrodes, patch clamp), behaviour testing (open field test, elevated plus maze, light-dark box test)"
$newSkills = ": PCR, qPCR, RT-PCR, gel electrophoresis, wire myography, western blotting, ELISA, microscopy, cell culture, immunohistochemistry, intracellular recording (microelectrodes, patch clamp), behaviour tests (open field test, elevated plus maze, light-dark box test)"
$skillsFound = $eduText.Find($oldSkills)
$skillsFound.Text = $newSkills

# --- Shape "Google Shape;101;p1": date label next to MSc entry ---
$msDate = $s.Shapes.Item("Google Shape;101;p1")
$msDate.Left = 4.689

# --- Shape "Google Shape;102;p1": date label next to BSc entry ---
$bscDate = $s.Shapes.Item("Google Shape;102;p1")
$bscDate.Left = 4.6902

# --- Shape "Google Shape;105;p1": personal statement paragraph ---
$quote = [char]0x2019
$statement = $s.Shapes.Item("Google Shape;105;p1")
$newStatement = "I am a curious and dedicated graduate with a Master" + $quote + "s degree in Medical Biology, driven by a genuine passion for advancing healthcare through science. My experience spans both fundamental research and clinical trial development. With a solid biomedical background and a keen interest for sharing knowledge, I am motivated to do meaningful research to support the health and well-being of both people and animals."
$statement.TextFrame.TextRange.Text = $newStatement
